$wb = $excel.ActiveWorkbook

# The "Swiss" sheet is a new market sheet, structurally identical to the
# existing "Czech" sheet (same layout/styles), so clone it and place the
# copy right after "Czech" (i.e. before "Slot Cards 215 Panel").
$czech = $wb.Worksheets.Item("Czech")
$czech.Copy($null, $czech)
$swiss = $wb.Worksheets.Item("Czech (2)")
$swiss.Name = "Swiss"

# Czech has an extra "FBI800" row (row 9) that Swiss doesn't use - drop it.
$swiss.Rows.Item(9).Delete()

# Row 8 keeps a custom row height on the new sheet.
$swiss.Rows.Item(8).RowHeight = 15

# Fill in the market-specific values.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2344"

# Czech is no longer the active tab; its selection resets to "select all".
$czech.Cells.Select() | Out-Null

# Swiss becomes the active/selected sheet and tab, with B2:B4 selected.
$swiss.Activate() | Out-Null
$swiss.Range("B2:B4").Select() | Out-Null
